$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the three new header cells (Wins / Losses / Ties) after the last
# existing column (AC). Copy the formatting of the neighboring header
# cell first so the new headers match the existing bold/bordered/centered
# look, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row on this sheet gets the team's overall W/L/T record.
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 79
    $ws.Cells.Item($r, 31).Value = 83
    $ws.Cells.Item($r, 32).Value = 0
}
